$d = $word.ActiveDocument

# 1. Update the document title in the header table:
#    "Circle Language Spec: Interfaces" -> "Circle Language Broader View | Interfaces"
$d.Content.Find.Execute(
    "Circle Language Spec: Interfaces", $true, $false, $false, $false, $false,
    $true, 1, $false, "Circle Language Broader View | Interfaces", 2) | Out-Null

# 2. Remove the stray "_GoBack" bookmark (an empty bookmarkStart/bookmarkEnd pair
#    Word automatically drops the next time it saves after an edit session).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Clean up the "eachother" spell-check split: re-typing/replacing the sentence
#    merges the three runs (and drops the proofErr spell-check markers) into one run.
$oldSentence = "an agreement with eachother, not to change the interface just like that. And at times it "
$newSentence = "an agreement with eachother, not to change the interface just like that. And at times it "
$d.Content.Find.Execute(
    $oldSentence, $true, $false, $false, $false, $false,
    $true, 1, $false, $newSentence, 2) | Out-Null
